$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.01259263455413108
$ws.Range("C2").Value = 0.6033495300249013
$ws.Range("D2").Value = 0.7009422445240612
$ws.Range("E2").Value = 0.8372229359758733
$ws.Range("F2").Value = 0.845458067251492
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = 0.2836421595116398
$ws.Range("C3").Value = 0.5962606503161496
$ws.Range("D3").Value = 0.6021311585356434
$ws.Range("E3").Value = 0.7759711067659951
$ws.Range("F3").Value = 0.7296059069611965
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = 0.1183466818900022
$ws.Range("C4").Value = 0.7836512297193149
$ws.Range("D4").Value = 1.193643475559592
$ws.Range("E4").Value = 1.092539919435254
$ws.Range("F4").Value = 1.097366538808324
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = 0.2525900972880621
$ws.Range("C5").Value = 0.6703098833254503
$ws.Range("D5").Value = 1.002878687496652
$ws.Range("E5").Value = 1.001438309381387
$ws.Range("F5").Value = 0.9793147045085741
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = 0.0217409311959274
$ws.Range("C6").Value = 0.7004675121762703
$ws.Range("D6").Value = 1.024571082396416
$ws.Range("E6").Value = 1.012210987095287
$ws.Range("F6").Value = 1.022918092182276
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = 0.1372560159180706
$ws.Range("C7").Value = 0.6144005222385287
$ws.Range("D7").Value = 0.9138507028447672
$ws.Range("E7").Value = 0.9559553874761977
$ws.Range("F7").Value = 0.9587497007249187
$ws.Range("G7").Value = 38

$ws.Range("B8").Value = 0.09871507743660377
$ws.Range("C8").Value = 0.6852141803533275
$ws.Range("D8").Value = 1.025442875956981
$ws.Range("E8").Value = 1.012641533790206
$ws.Range("F8").Value = 1.021720142012909
$ws.Range("G8").Value = 37

$ws.Range("B9").Value = 0.1617831974124649
$ws.Range("C9").Value = 0.7534609121388414
$ws.Range("D9").Value = 1.440524723083455
$ws.Range("E9").Value = 1.200218614704611
$ws.Range("F9").Value = 1.220160006814661
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -0.2070359202823368
$ws.Range("C10").Value = 0.6345706742736936
$ws.Range("D10").Value = 1.392524628435378
$ws.Range("E10").Value = 1.180052807477436
$ws.Range("F10").Value = 1.209186704288152
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = 0.02035458815990314
$ws.Range("C11").Value = 0.3973775210872706
$ws.Range("D11").Value = 0.2988378101058986
$ws.Range("E11").Value = 0.5466605986404166
$ws.Range("F11").Value = 0.6107613085800574
$ws.Range("G11").Value = 5

